$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing test case to the new iAU_TC_ID_212 entry
$ws.Range("A2").Value = " iAU_TC_ID_212"
$ws.Range("B2").Value = "@RegressionA Pre-Request Verify Elumina Login and Create Exam"
$ws.Range("C2").Value = "passed"

# Row 3: new Pre-Request row, reuses the same Test Case ID as row 2
$ws.Range("A3").Value = " iAU_TC_ID_212"
$ws.Range("B3").Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Range("C3").Value = "passed"

# Row 4: new iAU_TC_ID_218 entry
$ws.Range("A4").Value = "iAU_TC_ID_218"
$ws.Range("B4").Value = "@RegressionA Validation of Delivery--> Live Monitor - Candidate answer response Validation"
$ws.Range("C4").Value = "timedOut"

# Row 5: new iAU_TC_ID_219 entry
$ws.Range("A5").Value = "iAU_TC_ID_219"
$ws.Range("B5").Value = "@RegressionA Validation of Delivery--> Live Monitor - Live Streaming page"
$ws.Range("C5").Value = "failed"
